$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing IMF - Sales / IMF - Sales + Emp values (columns F and G)
# before they get overwritten, so we can move them into columns H and I
# (which previously held the now-removed "OECD (20%)" columns).
$oldF2 = $ws.Range("F2").Value2
$oldF3 = $ws.Range("F3").Value2
$oldF4 = $ws.Range("F4").Value2
$oldF5 = $ws.Range("F5").Value2
$oldF6 = $ws.Range("F6").Value2

$oldG2 = $ws.Range("G2").Value2
$oldG3 = $ws.Range("G3").Value2
$oldG4 = $ws.Range("G4").Value2
$oldG5 = $ws.Range("G5").Value2
$oldG6 = $ws.Range("G6").Value2

# Update header row: a new "IMF (20%)" pair of columns is inserted in the
# place of columns F/G, the old "IMF" columns move to H/I, replacing the
# "OECD (20%)" columns that are removed entirely.
$ws.Range("F1").Value = "IMF (20%) - Sales"
$ws.Range("G1").Value = "IMF (20%) - Sales + Emp"
$ws.Range("H1").Value = "IMF - Sales"
$ws.Range("I1").Value = "IMF - Sales + Emp"

# Move the previous "IMF - Sales"/"IMF - Sales + Emp" data into the H/I columns.
$ws.Range("H2").Value = $oldF2
$ws.Range("H3").Value = $oldF3
$ws.Range("H4").Value = $oldF4
$ws.Range("H5").Value = $oldF5
$ws.Range("H6").Value = $oldF6

$ws.Range("I2").Value = $oldG2
$ws.Range("I3").Value = $oldG3
$ws.Range("I4").Value = $oldG4
$ws.Range("I5").Value = $oldG5
$ws.Range("I6").Value = $oldG6

# Populate the new "IMF (20%)" data into columns F/G.
$ws.Range("F2").Value = 0.007444774198808549
$ws.Range("F3").Value = 0.02973125562628548
$ws.Range("F4").Value = -0.0005012849520288268
$ws.Range("F5").Value = -0.06426616309521044
$ws.Range("F6").Value = 0.002446805784062994

$ws.Range("G2").Value = 0.00606637845192245
$ws.Range("G3").Value = 0.0434634358097193
$ws.Range("G4").Value = 0.01716044213773076
$ws.Range("G5").Value = -0.08563732635657206
$ws.Range("G6").Value = 0.005593890638537687
